$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user changed which value the AutoFilter on column F (Project Name) is
# filtering for: from "State Regulatory Compliance" to "The Enlightenment Portal".
# Re-applying the AutoFilter causes Excel to recompute which rows are hidden
# (rows matching the new value become visible, others become hidden).
$xlFilterValues = 7
$ws.Range("F1:F140").AutoFilter(1, @("The Enlightenment Portal"), $xlFilterValues) | Out-Null

# Selection moved to G28 as a result of the user's subsequent interaction.
$ws.Range("G28").Select() | Out-Null
